# edit.ps1 - Apply the "fix: replace package with correct metadata file" edit.
#
# Summary of changes (from the OOXML diff):
#  1. Rename sheet "charts" -> "visualizations"
#  2. "Package info" sheet: bump DHIS2 version / Created timestamp / Identifier
#  3. "dashboardItems" sheet: replace generic "Chart"/"Map" content-type labels
#     with their specific visualization type (SINGLE_VALUE, PIE, COLUMN, BAR,
#     STACKED_COLUMN, STACKED_BAR, MAP); for MAP rows also blank out the
#     Content UID (col A) and Content name (col C).
#  4. "visualizations" (formerly "charts") sheet: blank out the stray
#     single-space Description cells (col B) that had no real description.
#  5. "programs" sheet: bump the "Last updated" date for the AEFI program.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the "charts" sheet to "visualizations"
# ---------------------------------------------------------------------------
$wsCharts = $wb.Worksheets.Item("charts")
$wsCharts.Name = "visualizations"

# ---------------------------------------------------------------------------
# 2. Package info sheet
# ---------------------------------------------------------------------------
$wsPkg = $wb.Worksheets.Item("Package info")
$wsPkg.Range("B5").Value = "DHIS2.34.4-aff07fb"
$wsPkg.Range("B6").Value = "20210406T141800"
$wsPkg.Range("B7").Value = "AEFI_TRACKER_V1.1.2_DHIS2.34.4-aff07fb_20210406T141800"

# ---------------------------------------------------------------------------
# 3. dashboardItems sheet: Content/item type (column B) corrections
# ---------------------------------------------------------------------------
$wsItems = $wb.Worksheets.Item("dashboardItems")

$itemTypes = @{
    3  = "SINGLE_VALUE"
    4  = "SINGLE_VALUE"
    5  = "MAP"
    6  = "PIE"
    7  = "PIE"
    8  = "SINGLE_VALUE"
    9  = "SINGLE_VALUE"
    10 = "PIE"
    11 = "PIE"
    12 = "PIE"
    13 = "STACKED_COLUMN"
    14 = "PIE"
    15 = "PIE"
    16 = "PIE"
    17 = "PIE"
    18 = "PIE"
    19 = "PIE"
    20 = "PIE"
    21 = "PIE"
    22 = "PIE"
    23 = "PIE"
    24 = "PIE"
    25 = "STACKED_BAR"
    26 = "BAR"
    27 = "COLUMN"
    30 = "COLUMN"
    31 = "COLUMN"
    32 = "COLUMN"
    33 = "COLUMN"
    34 = "MAP"
    35 = "MAP"
    36 = "MAP"
    37 = "SINGLE_VALUE"
    38 = "SINGLE_VALUE"
    39 = "SINGLE_VALUE"
    40 = "COLUMN"
    41 = "COLUMN"
    42 = "COLUMN"
    43 = "STACKED_BAR"
    44 = "STACKED_BAR"
    45 = "STACKED_BAR"
    47 = "PIE"
    48 = "PIE"
    49 = "PIE"
    50 = "PIE"
}

foreach ($r in $itemTypes.Keys) {
    $wsItems.Cells.Item($r, 2).Value = $itemTypes[$r]
}

# MAP rows additionally lose their Content UID (A) and Content name (C) -
# they now point at nothing in particular.
$mapRows = @(5, 34, 35, 36)
foreach ($r in $mapRows) {
    $wsItems.Cells.Item($r, 1).Value = ""
    $wsItems.Cells.Item($r, 3).Value = ""
}

# ---------------------------------------------------------------------------
# 4. visualizations (formerly charts) sheet: clear the placeholder
#    single-space Description cells (column B)
# ---------------------------------------------------------------------------
$wsViz = $wb.Worksheets.Item("visualizations")

$blankDescRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,30,31,32,36,37,38,39,40,41,42)
foreach ($r in $blankDescRows) {
    $wsViz.Cells.Item($r, 2).Value = ""
}

# ---------------------------------------------------------------------------
# 5. programs sheet: bump "Last updated" date
#    (force Text format first so the date-shaped string isn't silently
#    reinterpreted as a serial date - the source file keeps it as a plain
#    string, e.g. "2021-03-19" -> "2021-04-06")
# ---------------------------------------------------------------------------
$wsPrograms = $wb.Worksheets.Item("programs")
$wsPrograms.Range("C2").NumberFormat = "@"
$wsPrograms.Range("C2").Value = "2021-04-06"
